$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$name = "Nalathni Dragon"
$cost = "{2}{R}{R}"
$type = "Creature — Dragon"
$text = "Flying; banding (Any creatures with banding, and up to one without, can attack in a band. Bands are blocked as a group. If any creatures with banding you control are blocking or being blocked by a creature, you divide that creature’s combat damage, not its controller, among any of the creatures it’s being blocked by or is blocking.)"
$ability = "{R}: Nalathni Dragon gets +1/+0 until end of turn. If this ability has been activated four or more times this turn, sacrifice Nalathni Dragon at the beginning of the next end step."
$pt = "1/1"

$combined = "('" + $name + "', ['" + $cost + "', '" + $type + "', '" + $text + "', '" + $ability + "', '" + $pt + "'])"

$ws.Range("A2").Value = $combined

$ws.Range("A3:A7").EntireRow.Delete()
